# Applies the "Renames" sheet addition to the workbook.
$wb = $excel.ActiveWorkbook

# --- Add the new "Renames" worksheet at the end (after the last sheet) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Renames"

# --- Populate the Renames sheet with header + data rows ---
# (column B entered before column A so new shared-string entries land in the
# same order as the authored workbook: "new name" then "old name")
$ws.Range("B1").Value = "new name"
$ws.Range("A1").Value = "old name"

$ws.Range("B2").Value = "Big Cylinder"
$ws.Range("A2").Value = "Cylinder"

$ws.Range("B3").Value = "Locking Nut"
$ws.Range("A3").Value = "Lug Nut"

# --- Column widths / best-fit like the source sheet ---
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null

# --- Select A4 on the new sheet (matches target selection) ---
$ws.Range("A4").Select()
